$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 86.360967
$ws.Range("H2").Value = 259.082901
$ws.Range("I2").Value = 0.2482072712525276
$ws.Range("J2").Value = 0.2482072712525276
$ws.Range("M2").Value = 255.0443116666667
$ws.Range("N2").Value = 765.132935
$ws.Range("O2").Value = 0.863617428561108
$ws.Range("P2").Value = 0.8636174285611079
$ws.Range("Q2").Value = 22025.87338338271
$ws.Range("R2").Value = 198232.8604504444
$ws.Range("S2").Value = 0.2143561253492773
$ws.Range("T2").Value = 0.2143561253492773

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 86.360967
$ws.Range("H3").Value = 259.082901
$ws.Range("I3").Value = 0.2482072712525276
$ws.Range("J3").Value = 0.2482072712525276
$ws.Range("M3").Value = 0.8952453333333334
$ws.Range("O3").Value = 0.003031431940796009
$ws.Range("P3").Value = 0.003031431940796009
$ws.Range("Q3").Value = 77.31425268890401
$ws.Range("R3").Value = 695.8282742001361
$ws.Range("S3").Value = 0.0007524234500127311
$ws.Range("T3").Value = 0.0007524234500127312

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 86.360967
$ws.Range("H4").Value = 259.082901
$ws.Range("I4").Value = 0.2482072712525276
$ws.Range("J4").Value = 0.2482072712525276
$ws.Range("M4").Value = 7.050555333333333
$ws.Range("N4").Value = 21.151666
$ws.Range("O4").Value = 0.02387421396349043
$ws.Range("P4").Value = 0.02387421396349043
$ws.Range("Q4").Value = 608.892776473674
$ws.Range("R4").Value = 5480.034988263065
$ws.Range("S4").Value = 0.005925753501176951
$ws.Range("T4").Value = 0.005925753501176951

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 86.360967
$ws.Range("H5").Value = 259.082901
$ws.Range("I5").Value = 0.2482072712525276
$ws.Range("J5").Value = 0.2482072712525276
$ws.Range("M5").Value = 32.33082866666666
$ws.Range("N5").Value = 96.99248599999999
$ws.Range("O5").Value = 0.1094769255346056
$ws.Range("P5").Value = 0.1094769255346056
$ws.Range("Q5").Value = 2792.121627564654
$ws.Range("R5").Value = 25129.09464808188
$ws.Range("S5").Value = 0.02717296895206063
$ws.Range("T5").Value = 0.02717296895206063

$ws.Range("I6").Value = 0.6003523616657895
$ws.Range("J6").Value = 0.6003523616657896
$ws.Range("M6").Value = 255.0443116666667
$ws.Range("N6").Value = 765.132935
$ws.Range("O6").Value = 0.863617428561108
$ws.Range("P6").Value = 0.8636174285611079
$ws.Range("Q6").Value = 53275.17214438097
$ws.Range("R6").Value = 479476.5492994287
$ws.Range("S6").Value = 0.5184747628123975
$ws.Range("T6").Value = 0.5184747628123975

$ws.Range("I7").Value = 0.6003523616657895
$ws.Range("J7").Value = 0.6003523616657896
$ws.Range("M7").Value = 0.8952453333333334
$ws.Range("O7").Value = 0.003031431940796009
$ws.Range("P7").Value = 0.003031431940796009
$ws.Range("S7").Value = 0.001819927324885992
$ws.Range("T7").Value = 0.001819927324885992

$ws.Range("I8").Value = 0.6003523616657895
$ws.Range("J8").Value = 0.6003523616657896
$ws.Range("M8").Value = 7.050555333333333
$ws.Range("N8").Value = 21.151666
$ws.Range("O8").Value = 0.02387421396349043
$ws.Range("P8").Value = 0.02387421396349043
$ws.Range("Q8").Value = 1472.761915928308
$ws.Range("R8").Value = 13254.85724335478
$ws.Range("S8").Value = 0.01433294073589585
$ws.Range("T8").Value = 0.01433294073589585

$ws.Range("I9").Value = 0.6003523616657895
$ws.Range("J9").Value = 0.6003523616657896
$ws.Range("M9").Value = 32.33082866666666
$ws.Range("N9").Value = 96.99248599999999
$ws.Range("O9").Value = 0.1094769255346056
$ws.Range("P9").Value = 0.1094769255346056
$ws.Range("Q9").Value = 6753.455709446699
$ws.Range("R9").Value = 60781.10138502029
$ws.Range("S9").Value = 0.06572473079261028
$ws.Range("T9").Value = 0.06572473079261028

$ws.Range("G10").Value = 52.26262533333333
$ws.Range("H10").Value = 156.787876
$ws.Range("I10").Value = 0.1502063266901572
$ws.Range("J10").Value = 0.1502063266901572
$ws.Range("M10").Value = 255.0443116666667
$ws.Range("N10").Value = 765.132935
$ws.Range("O10").Value = 0.863617428561108
$ws.Range("P10").Value = 0.8636174285611079
$ws.Range("Q10").Value = 13329.28530403289
$ws.Range("R10").Value = 119963.567736296
$ws.Range("S10").Value = 0.1297208016097633
$ws.Range("T10").Value = 0.1297208016097633

$ws.Range("G11").Value = 52.26262533333333
$ws.Range("H11").Value = 156.787876
$ws.Range("I11").Value = 0.1502063266901572
$ws.Range("J11").Value = 0.1502063266901572
$ws.Range("M11").Value = 0.8952453333333334
$ws.Range("O11").Value = 0.003031431940796009
$ws.Range("P11").Value = 0.003031431940796009
$ws.Range("Q11").Value = 46.78787143741511
$ws.Range("R11").Value = 421.090842936736
$ws.Range("S11").Value = 0.0004553402564381825
$ws.Range("T11").Value = 0.0004553402564381826

$ws.Range("G12").Value = 52.26262533333333
$ws.Range("H12").Value = 156.787876
$ws.Range("I12").Value = 0.1502063266901572
$ws.Range("J12").Value = 0.1502063266901572
$ws.Range("M12").Value = 7.050555333333333
$ws.Range("N12").Value = 21.151666
$ws.Range("O12").Value = 0.02387421396349043
$ws.Range("P12").Value = 0.02387421396349043
$ws.Range("Q12").Value = 368.480531777935
$ws.Range("R12").Value = 3316.324786001415
$ws.Range("S12").Value = 0.003586057982070756
$ws.Range("T12").Value = 0.003586057982070756

$ws.Range("G13").Value = 52.26262533333333
$ws.Range("H13").Value = 156.787876
$ws.Range("I13").Value = 0.1502063266901572
$ws.Range("J13").Value = 0.1502063266901572
$ws.Range("M13").Value = 32.33082866666666
$ws.Range("N13").Value = 96.99248599999999
$ws.Range("O13").Value = 0.1094769255346056
$ws.Range("P13").Value = 0.1094769255346056
$ws.Range("Q13").Value = 1689.693985322192
$ws.Range("R13").Value = 15207.24586789973
$ws.Range("S13").Value = 0.01644412684188499
$ws.Range("T13").Value = 0.01644412684188499

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4293706666666666
$ws.Range("H14").Value = 1.288112
$ws.Range("I14").Value = 0.001234040391525629
$ws.Range("J14").Value = 0.001234040391525629
$ws.Range("M14").Value = 255.0443116666667
$ws.Range("N14").Value = 765.132935
$ws.Range("O14").Value = 0.863617428561108
$ws.Range("P14").Value = 0.8636174285611079
$ws.Range("Q14").Value = 109.5085461298578
$ws.Range("R14").Value = 985.5769151687199
$ws.Range("S14").Value = 0.001065738789669907
$ws.Range("T14").Value = 0.001065738789669907

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4293706666666666
$ws.Range("H15").Value = 1.288112
$ws.Range("I15").Value = 0.001234040391525629
$ws.Range("J15").Value = 0.001234040391525629
$ws.Range("M15").Value = 0.8952453333333334
$ws.Range("O15").Value = 0.003031431940796009
$ws.Range("P15").Value = 0.003031431940796009
$ws.Range("Q15").Value = 0.3843920856035556
$ws.Range("R15").Value = 3.459528770432
$ws.Range("S15").Value = 0.000003740909459103204
$ws.Range("T15").Value = 0.000003740909459103205

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4293706666666666
$ws.Range("H16").Value = 1.288112
$ws.Range("I16").Value = 0.001234040391525629
$ws.Range("J16").Value = 0.001234040391525629
$ws.Range("M16").Value = 7.050555333333333
$ws.Range("N16").Value = 21.151666
$ws.Range("O16").Value = 0.02387421396349043
$ws.Range("P16").Value = 0.02387421396349043
$ws.Range("Q16").Value = 3.027301643843555
$ws.Range("R16").Value = 27.245714794592
$ws.Range("S16").Value = 0.00002946174434687237
$ws.Range("T16").Value = 0.00002946174434687237

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4293706666666666
$ws.Range("H17").Value = 1.288112
$ws.Range("I17").Value = 0.001234040391525629
$ws.Range("J17").Value = 0.001234040391525629
$ws.Range("M17").Value = 32.33082866666666
$ws.Range("N17").Value = 96.99248599999999
$ws.Range("O17").Value = 0.1094769255346056
$ws.Range("P17").Value = 0.1094769255346056
$ws.Range("Q17").Value = 13.88190945849244
$ws.Range("R17").Value = 124.937185126432
$ws.Range("S17").Value = 0.0001350989480497469
$ws.Range("T17").Value = 0.0001350989480497469
